$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the G and H column formulas for rows 3-91: the correlation weights flip
# from negative coefficients with a negative divisor to positive coefficients with
# a (near-)unity divisor.
$ws.Range("G3:G91").Formula = "=((0.078)*B3+(0.218)*C3+(0.477)*E3+(0.227)*F3)/(1)"
$ws.Range("H3:H91").Formula = "=((0.078)*B3+(0.207)*D3+(0.477)*E3+(0.227)*F3)/ (0.989)"

# Move/extend the selection to the newly (re)plotted centrality column G3:G91,
# anchored at G3.
$ws.Range("G3:G91").Select()
